$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 10, shifting the existing rows 10-15 down to 12-17.
$ws.Range("A10:T11").Insert()

# Row 10: new weekly entry (Primera)
$ws.Cells.Item(10, 1).Value = 1
$ws.Cells.Item(10, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(10, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(10, 4).Value = 44846
$ws.Cells.Item(10, 5).Value = 15
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100107
$ws.Cells.Item(10, 8).Value = "Otros"
$ws.Cells.Item(10, 9).Value = 100107002
$ws.Cells.Item(10, 10).Value = "Chirimoya"
$ws.Cells.Item(10, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 160
$ws.Cells.Item(10, 14).Value = 24000
$ws.Cells.Item(10, 15).Value = 25000
$ws.Cells.Item(10, 16).Value = 24500
$ws.Cells.Item(10, 17).Value = "$/caja 12 kilos"
$ws.Cells.Item(10, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(10, 19).Value = 2042
$ws.Cells.Item(10, 20).Value = 12

# Row 11: new weekly entry (Segunda)
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value = 44846
$ws.Cells.Item(11, 5).Value = 15
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100107
$ws.Cells.Item(11, 8).Value = "Otros"
$ws.Cells.Item(11, 9).Value = 100107002
$ws.Cells.Item(11, 10).Value = "Chirimoya"
$ws.Cells.Item(11, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(11, 12).Value = "Segunda"
$ws.Cells.Item(11, 13).Value = 100
$ws.Cells.Item(11, 14).Value = 22000
$ws.Cells.Item(11, 15).Value = 23000
$ws.Cells.Item(11, 16).Value = 22500
$ws.Cells.Item(11, 17).Value = "$/caja 12 kilos"
$ws.Cells.Item(11, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(11, 19).Value = 1875
$ws.Cells.Item(11, 20).Value = 12

# Apply the same date number-format (style index 2 in before.xlsx) to the new D cells.
$ws.Cells.Item(10, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
